$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.232.13'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '1.556.28'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.63'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3794'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.91%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3278'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.61'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -8.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.138'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07373'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.21'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.852'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.50%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.758'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.89%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.557.56'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001076'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06652'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '86.34'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.432'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.19'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.72'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('D24').Value = '22.246.68'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.291'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.567'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.58'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.945'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.93'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('D31').Value = '1.733.66'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.079'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.926'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.923'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.405'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08237'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02357'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06347'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.366'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2160'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.240'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.04'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6080'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.85'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.750'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5900'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.19'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.975'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.180'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07077'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.87%  '
